$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 (new response #29) ---
$ws.Range("A25").Value = 29
$ws.Range("B25").Value = "2024-07-04 18:08:19"
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = "de"
$ws.Range("E25").Value = 1147875795
$ws.Range("F25").Value = "2024-07-04 18:06:51"
$ws.Range("G25").Value = "2024-07-04 18:08:19"
$ws.Range("I25").Value = "TZ09CH25"
$ws.Range("J25").Value = "Ja"
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 3
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 1
$ws.Range("P25").Value = 3
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 3
$ws.Range("U25").Value = 3
$ws.Range("V25").Value = 3
$ws.Range("W25").Value = 3
$ws.Range("X25").Value = 3
$ws.Range("Y25").Value = 3
$ws.Range("Z25").Value = 4
$ws.Range("AA25").Value = 4
$ws.Range("AB25").Value = 3
$ws.Range("AC25").Value = 4
$ws.Range("AD25").Value = 3
$ws.Range("AE25").Value = 48
$ws.Range("AF25").Value = "Männlich"
$ws.Range("AH25").Value = "Promotion"
$ws.Range("AJ25").Value = 89.11
$ws.Range("AK25").Value = 14.7
$ws.Range("AN25").Value = 65.06
$ws.Range("AP25").Value = 9.35

# --- Row 26 (new response #30) ---
$ws.Range("A26").Value = 30
$ws.Range("B26").Value = "2024-07-06 12:34:45"
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = "de"
$ws.Range("E26").Value = 2043447202
$ws.Range("F26").Value = "2024-07-06 12:31:14"
$ws.Range("G26").Value = "2024-07-06 12:34:45"
$ws.Range("I26").Value = "ER09AS09"
$ws.Range("J26").Value = "Ja"
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 2
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 1
$ws.Range("P26").Value = 1
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 3
$ws.Range("T26").Value = 1
$ws.Range("U26").Value = 1
$ws.Range("V26").Value = 1
$ws.Range("W26").Value = 3
$ws.Range("X26").Value = 1
$ws.Range("Y26").Value = 2
$ws.Range("Z26").Value = 3
$ws.Range("AA26").Value = 4
$ws.Range("AB26").Value = 3
$ws.Range("AC26").Value = 2
$ws.Range("AD26").Value = 3
$ws.Range("AE26").Value = 28
$ws.Range("AF26").Value = "Männlich"
$ws.Range("AH26").Value = "Master-Abschluss"
$ws.Range("AJ26").Value = 211.78
$ws.Range("AK26").Value = 104.08
$ws.Range("AN26").Value = 77.58
$ws.Range("AP26").Value = 30.12

# Update the view: select E26 (the cell the author last edited) and
# scroll the sheet back to the default top-left position.
$ws.Range("A1").Select() | Out-Null
$ws.Range("E26").Select() | Out-Null
